# Add problem_2 and problem_3 according to problem suite
#
# 1. Sheet1: update the "Row:" / "Col:" parameters used by the generator
#    (B1: 6 -> 11, B2: 6 -> 7) and move the selection to C1.
# 2. problem: replace the 6x6 sample maze with the newly generated
#    9-wide x 14-tall maze (problem_2 / problem_3) and move the selection
#    to I1.
# 3. Workbook: register the new defined name "A" pointing at problem!$XDD$9
#    (left behind by the generator macro run).

$wb = $excel.ActiveWorkbook

# --- Sheet1: generator parameters -----------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B1").Value = 11
$ws1.Range("B2").Value = 7
$ws1.Range("C1").Select() | Out-Null

# --- problem: new maze grid (A1:I14) ---------------------------------------
$ws2 = $wb.Worksheets.Item("problem")

$grid = @(
    @(3,3,3,1,1,1,3,3,3),
    @(3,3,3,1,4,1,3,3,3),
    @(3,1,1,1,2,1,1,1,3),
    @(1,1,0,0,5,0,0,1,1),
    @(1,0,0,1,0,1,0,0,1),
    @(1,0,1,0,0,0,1,0,1),
    @(1,0,1,0,0,0,1,0,1),
    @(1,0,1,0,0,0,1,0,1),
    @(1,0,0,1,0,1,0,0,1),
    @(1,1,0,2,0,2,0,1,1),
    @(3,1,1,5,0,5,1,1,3),
    @(3,3,1,0,0,0,1,3,3),
    @(3,3,1,0,0,0,1,3,3),
    @(3,3,1,1,1,1,1,3,3)
)

for ($r = 0; $r -lt $grid.Length; $r++) {
    $row = $grid[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws2.Range("I1").Select() | Out-Null

# --- Workbook: new defined name --------------------------------------------
$wb.Names.Add('A', '=problem!$XDD$9')
